$p = $ppt.ActivePresentation

# --- 1) Swap the deck's theme colour scheme from "Integral" (Red Violet)
#        to the stock "Office Theme" (Office) palette. The ThemeColorScheme
#        exposed on a slide edits the presentation's underlying theme part
#        (ppt/theme/theme1.xml) in place - set every one of the twelve
#        standard theme colour slots to its "Office" RGB value (values use
#        the Windows BGR-packed COLORREF order expected by the RGB setter).
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

$tcs.Item(1).RGB  = 0         # Dark 1    -> 000000
$tcs.Item(2).RGB  = 16777215  # Light 1   -> FFFFFF
$tcs.Item(3).RGB  = 6968388   # Dark 2    -> 44546A
$tcs.Item(4).RGB  = 15132391  # Light 2   -> E7E6E6
$tcs.Item(5).RGB  = 13998939  # Accent 1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501   # Accent 2  -> ED7D31
$tcs.Item(7).RGB  = 10855845  # Accent 3  -> A5A5A5
$tcs.Item(8).RGB  = 49407     # Accent 4  -> FFC000
$tcs.Item(9).RGB  = 12874308  # Accent 5  -> 4472C4
$tcs.Item(10).RGB = 4697456   # Accent 6  -> 70AD47
$tcs.Item(11).RGB = 12673797  # Hyperlink -> 0563C1
$tcs.Item(12).RGB = 7491477   # Followed Hyperlink -> 954F72

# --- 2) Re-style the three tables (slides 14, 15 and 16) away from the
#        deck's custom "Table_0" style onto the built-in style the author
#        picked from the Table Design gallery.
$newTableStyle = "{BB563EDC-6DB7-41E0-A7DE-3CE57211825A}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}
